$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A180").Value = "IMX-USD"
$ws.Range("A181").Value = "TAO-USD"
$ws.Range("A182").Value = "GRT-USD"
